# Auto-generated edit script applying the Behemoth_Profits.xlsx diff
# Updates market-price columns (H:N) across all 8 sheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Cells.Item(8,8).Value = 2498.4285  # H8: 1596.2727 -> 2498.4285
$ws.Cells.Item(8,9).Value = 2498.4285  # I8: 1754.7 -> 2498.4285
$ws.Cells.Item(8,10).Value = 0  # J8: 12 -> 0
$ws.Cells.Item(8,11).Value = 7495.2855  # K8: 5264.1 -> 7495.2855
$ws.Cells.Item(8,12).Value = 0  # L8: 36 -> 0
$ws.Cells.Item(8,13).Value = -7356.2855  # M8: -5125.1 -> -7356.2855
$ws.Cells.Item(8,14).ClearContents()  # N8 removed (was -314)

# Row 19
$ws.Cells.Item(19,8).Value = 2157.2856  # H19: 1810.0588 -> 2157.2856
$ws.Cells.Item(19,9).Value = 1740.5  # I19: 1460 -> 1740.5
$ws.Cells.Item(19,10).Value = 2469.875  # J19: 2203.875 -> 2469.875
$ws.Cells.Item(19,11).Value = 1740.5  # K19: 1460 -> 1740.5
$ws.Cells.Item(19,12).Value = 2469.875  # L19: 2203.875 -> 2469.875
$ws.Cells.Item(19,13).Value = -1565.5  # M19: -1285 -> -1565.5
$ws.Cells.Item(19,14).Value = -2819.875  # N19: -2553.875 -> -2819.875

# Row 64
$ws.Cells.Item(64,8).Value = 4999.9414  # H64: 4999.9473 -> 4999.9414
$ws.Cells.Item(64,10).Value = 4999.9414  # J64: 4999.9473 -> 4999.9414
$ws.Cells.Item(64,12).Value = 4999.9414  # L64: 4999.9473 -> 4999.9414
$ws.Cells.Item(64,14).Value = -5495.9414  # N64: -5495.9473 -> -5495.9414

# Row 67
$ws.Cells.Item(67,8).Value = 4999.9414  # H67: 4999.9473 -> 4999.9414
$ws.Cells.Item(67,10).Value = 4999.9414  # J67: 4999.9473 -> 4999.9414
$ws.Cells.Item(67,12).Value = 4999.9414  # L67: 4999.9473 -> 4999.9414
$ws.Cells.Item(67,14).Value = -6715.9414  # N67: -6715.9473 -> -6715.9414

# Row 70
$ws.Cells.Item(70,8).Value = 1494.4166  # H70: 1520.3478 -> 1494.4166
$ws.Cells.Item(70,9).Value = 1630.2727  # I70: 1703.5 -> 1630.2727
$ws.Cells.Item(70,11).Value = 4890.8181  # K70: 5110.5 -> 4890.8181
$ws.Cells.Item(70,13).Value = -4620.8181  # M70: -4840.5 -> -4620.8181

# Row 73
$ws.Cells.Item(73,8).Value = 1494.4166  # H73: 1520.3478 -> 1494.4166
$ws.Cells.Item(73,9).Value = 1630.2727  # I73: 1703.5 -> 1630.2727
$ws.Cells.Item(73,11).Value = 4890.8181  # K73: 5110.5 -> 4890.8181
$ws.Cells.Item(73,13).Value = -3954.8181  # M73: -4174.5 -> -3954.8181

# Row 101
$ws.Cells.Item(101,8).Value = 2234.2727  # H101: 2352.4546 -> 2234.2727
$ws.Cells.Item(101,9).Value = 2286.3333  # I101: 2447.125 -> 2286.3333
$ws.Cells.Item(101,10).Value = 2000  # J101: 2100 -> 2000
$ws.Cells.Item(101,11).Value = 6858.999899999999  # K101: 7341.375 -> 6858.999899999999
$ws.Cells.Item(101,12).Value = 6000  # L101: 6300 -> 6000
$ws.Cells.Item(101,13).Value = -5236.999899999999  # M101: -5719.375 -> -5236.999899999999
$ws.Cells.Item(101,14).Value = -9244  # N101: -9544 -> -9244

# Row 109
$ws.Cells.Item(109,8).Value = 115990  # H109: 116000 -> 115990
$ws.Cells.Item(109,10).Value = 115990  # J109: 116000 -> 115990
$ws.Cells.Item(109,12).Value = 115990  # L109: 116000 -> 115990
$ws.Cells.Item(109,14).Value = -118764  # N109: -118774 -> -118764

# Row 110
$ws.Cells.Item(110,8).Value = 33248.5  # H110: 33248.75 -> 33248.5
$ws.Cells.Item(110,10).Value = 33248.5  # J110: 33248.75 -> 33248.5
$ws.Cells.Item(110,12).Value = 33248.5  # L110: 33248.75 -> 33248.5
$ws.Cells.Item(110,14).Value = -41428.5  # N110: -41428.75 -> -41428.5

# Row 116
$ws.Cells.Item(116,8).Value = 7017.5454  # H116: 6641.0835 -> 7017.5454
$ws.Cells.Item(116,9).Value = 6288.125  # I116: 5867.222 -> 6288.125
$ws.Cells.Item(116,11).Value = 6288.125  # K116: 5867.222 -> 6288.125
$ws.Cells.Item(116,13).Value = -2846.125  # M116: -2425.222 -> -2846.125

# Row 126
$ws.Cells.Item(126,8).Value = 72110  # H126: 73447.5 -> 72110
$ws.Cells.Item(126,10).Value = 67775  # J126: 68790 -> 67775
$ws.Cells.Item(126,12).Value = 67775  # L126: 68790 -> 67775
$ws.Cells.Item(126,14).Value = -77655  # N126: -78670 -> -77655

# Row 138
$ws.Cells.Item(138,8).Value = 3059.0725  # H138: 3077.913 -> 3059.0725
$ws.Cells.Item(138,10).Value = 3173.554  # J138: 3193.554 -> 3173.554
$ws.Cells.Item(138,12).Value = 9520.662  # L138: 9580.662 -> 9520.662
$ws.Cells.Item(138,14).Value = -19800.662  # N138: -19860.662 -> -19800.662

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45,8).Value = 2588.8667  # H45: 2526.3845 -> 2588.8667
$ws.Cells.Item(45,9).Value = 1932  # I45: 1925.5555 -> 1932
$ws.Cells.Item(45,10).Value = 3902.6  # J45: 3878.25 -> 3902.6
$ws.Cells.Item(45,11).Value = 1932  # K45: 1925.5555 -> 1932
$ws.Cells.Item(45,12).Value = 3902.6  # L45: 3878.25 -> 3902.6
$ws.Cells.Item(45,13).Value = -1555  # M45: -1548.5555 -> -1555
$ws.Cells.Item(45,14).Value = -4656.6  # N45: -4632.25 -> -4656.6

# Row 103
$ws.Cells.Item(103,8).Value = 54570.4  # H103: 55838 -> 54570.4
$ws.Cells.Item(103,10).Value = 54570.4  # J103: 55838 -> 54570.4
$ws.Cells.Item(103,12).Value = 54570.4  # L103: 55838 -> 54570.4
$ws.Cells.Item(103,14).Value = -56914.4  # N103: -58182 -> -56914.4

# Row 105
$ws.Cells.Item(105,8).Value = 70365  # H105: 70370 -> 70365
$ws.Cells.Item(105,10).Value = 70365  # J105: 70370 -> 70365
$ws.Cells.Item(105,12).Value = 70365  # L105: 70370 -> 70365
$ws.Cells.Item(105,14).Value = -77353  # N105: -77358 -> -77353

# Row 122
$ws.Cells.Item(122,8).Value = 3276.2222  # H122: 2425.8572 -> 3276.2222
$ws.Cells.Item(122,9).Value = 1121.5  # I122: 870.25 -> 1121.5
$ws.Cells.Item(122,10).Value = 5000  # J122: 4500 -> 5000
$ws.Cells.Item(122,11).Value = 3364.5  # K122: 2610.75 -> 3364.5
$ws.Cells.Item(122,12).Value = 15000  # L122: 13500 -> 15000
$ws.Cells.Item(122,13).Value = -914.5  # M122: -160.75 -> -914.5
$ws.Cells.Item(122,14).Value = -19900  # N122: -18400 -> -19900

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20,8).Value = 3138.1  # H20: 3107.762 -> 3138.1
$ws.Cells.Item(20,9).Value = 3250.8823  # I20: 3251 -> 3250.8823
$ws.Cells.Item(20,11).Value = 3250.8823  # K20: 3251 -> 3250.8823
$ws.Cells.Item(20,13).Value = -3003.8823  # M20: -3004 -> -3003.8823

# Row 29
$ws.Cells.Item(29,8).Value = 3750  # H29: 0 -> 3750
$ws.Cells.Item(29,9).Value = 5000  # I29: 0 -> 5000
$ws.Cells.Item(29,10).Value = 2500  # J29: 0 -> 2500
$ws.Cells.Item(29,11).Value = 5000  # K29: 0 -> 5000
$ws.Cells.Item(29,12).Value = 2500  # L29: 0 -> 2500
$ws.Cells.Item(29,13).Value = -4711  # M29: None -> -4711
$ws.Cells.Item(29,14).Value = -3078  # N29: None -> -3078

# Row 102
$ws.Cells.Item(102,8).Value = 76429.2  # H102: 76469 -> 76429.2
$ws.Cells.Item(102,10).Value = 99863.336  # J102: 99929.664 -> 99863.336
$ws.Cells.Item(102,12).Value = 99863.336  # L102: 99929.664 -> 99863.336
$ws.Cells.Item(102,14).Value = -106353.336  # N102: -106419.664 -> -106353.336

# Row 105
$ws.Cells.Item(105,8).Value = 2277.9  # H105: 1499.3334 -> 2277.9
$ws.Cells.Item(105,9).Value = 1194.75  # I105: 899 -> 1194.75
$ws.Cells.Item(105,10).Value = 3000  # J105: 2700 -> 3000
$ws.Cells.Item(105,11).Value = 1194.75  # K105: 899 -> 1194.75
$ws.Cells.Item(105,12).Value = 3000  # L105: 2700 -> 3000
$ws.Cells.Item(105,13).Value = 552.25  # M105: 848 -> 552.25
$ws.Cells.Item(105,14).Value = -6494  # N105: -6194 -> -6494

# Row 107
$ws.Cells.Item(107,8).Value = 2912.6667  # H107: 3189.25 -> 2912.6667
$ws.Cells.Item(107,10).Value = 4625.75  # J107: 5934.3335 -> 4625.75
$ws.Cells.Item(107,12).Value = 4625.75  # L107: 5934.3335 -> 4625.75
$ws.Cells.Item(107,14).Value = -8465.75  # N107: -9774.333500000001 -> -8465.75

# Row 134
$ws.Cells.Item(134,8).Value = 76706.28999999999  # H134: 63457.117 -> 76706.28999999999
$ws.Cells.Item(134,9).Value = 1260  # I134: 1338.7858 -> 1260
$ws.Cells.Item(134,11).Value = 3780  # K134: 4016.3574 -> 3780
$ws.Cells.Item(134,13).Value = -1245  # M134: -1481.3574 -> -1245

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Cells.Item(58,8).Value = 1099.3334  # H58: 1139.2 -> 1099.3334
$ws.Cells.Item(58,10).Value = 1300  # J58: 1700 -> 1300
$ws.Cells.Item(58,12).Value = 1300  # L58: 1700 -> 1300
$ws.Cells.Item(58,14).Value = -1706  # N58: -2106 -> -1706

# Row 124
$ws.Cells.Item(124,8).Value = 0  # H124: 24980 -> 0
$ws.Cells.Item(124,10).Value = 0  # J124: 24980 -> 0
$ws.Cells.Item(124,12).Value = 0  # L124: 24980 -> 0
$ws.Cells.Item(124,14).ClearContents()  # N124 removed (was -29890)

# Row 136
$ws.Cells.Item(136,8).Value = 1099.3334  # H136: 1139.2 -> 1099.3334
$ws.Cells.Item(136,10).Value = 1300  # J136: 1700 -> 1300
$ws.Cells.Item(136,12).Value = 3900  # L136: 5100 -> 3900
$ws.Cells.Item(136,14).Value = -9000  # N136: -10200 -> -9000

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5,8).Value = 700  # H5: 724 -> 700
$ws.Cells.Item(5,10).Value = 900  # J5: 780 -> 900
$ws.Cells.Item(5,12).Value = 2700  # L5: 2340 -> 2700
$ws.Cells.Item(5,14).Value = -2924  # N5: -2564 -> -2924

# Row 29
$ws.Cells.Item(29,8).Value = 1094.8  # H29: 1095 -> 1094.8
$ws.Cells.Item(29,9).Value = 195  # I29: 0 -> 195
$ws.Cells.Item(29,10).Value = 1319.75  # J29: 1095 -> 1319.75
$ws.Cells.Item(29,11).Value = 585  # K29: 0 -> 585
$ws.Cells.Item(29,12).Value = 3959.25  # L29: 3285 -> 3959.25
$ws.Cells.Item(29,13).Value = -308  # M29: None -> -308
$ws.Cells.Item(29,14).Value = -4513.25  # N29: -3839 -> -4513.25

# Row 40
$ws.Cells.Item(40,8).Value = 270.14285  # H40: 252.46666 -> 270.14285
$ws.Cells.Item(40,9).Value = 217.2  # I40: 197.90909 -> 217.2
$ws.Cells.Item(40,11).Value = 868.8  # K40: 791.63636 -> 868.8
$ws.Cells.Item(40,13).Value = -799.8  # M40: -722.63636 -> -799.8

# Row 58
$ws.Cells.Item(58,8).Value = 1933  # H58: 1824.5 -> 1933
$ws.Cells.Item(58,10).Value = 1499.5  # J58: 1499.3334 -> 1499.5
$ws.Cells.Item(58,12).Value = 4498.5  # L58: 4498.0002 -> 4498.5
$ws.Cells.Item(58,14).Value = -4754.5  # N58: -4754.0002 -> -4754.5

# Row 135
$ws.Cells.Item(135,8).Value = 700  # H135: 724 -> 700
$ws.Cells.Item(135,10).Value = 900  # J135: 780 -> 900
$ws.Cells.Item(135,12).Value = 8100  # L135: 7020 -> 8100
$ws.Cells.Item(135,14).Value = -13170  # N135: -12090 -> -13170

# Row 137
$ws.Cells.Item(137,8).Value = 10000  # H137: 3925 -> 10000
$ws.Cells.Item(137,9).Value = 10000  # I137: 5950 -> 10000
$ws.Cells.Item(137,10).Value = 0  # J137: 1900 -> 0
$ws.Cells.Item(137,11).Value = 30000  # K137: 17850 -> 30000
$ws.Cells.Item(137,12).Value = 0  # L137: 5700 -> 0
$ws.Cells.Item(137,13).Value = -24900  # M137: -12750 -> -24900
$ws.Cells.Item(137,14).ClearContents()  # N137 removed (was -15900)

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70,8).Value = 5059.9  # H70: 4795.769 -> 5059.9
$ws.Cells.Item(70,9).Value = 4955.4443  # I70: 4695.4165 -> 4955.4443
$ws.Cells.Item(70,11).Value = 4955.4443  # K70: 4695.4165 -> 4955.4443
$ws.Cells.Item(70,13).Value = -4685.4443  # M70: -4425.4165 -> -4685.4443

# Row 73
$ws.Cells.Item(73,8).Value = 5059.9  # H73: 4795.769 -> 5059.9
$ws.Cells.Item(73,9).Value = 4955.4443  # I73: 4695.4165 -> 4955.4443
$ws.Cells.Item(73,11).Value = 4955.4443  # K73: 4695.4165 -> 4955.4443
$ws.Cells.Item(73,13).Value = -4019.4443  # M73: -3759.4165 -> -4019.4443

# Row 97
$ws.Cells.Item(97,8).Value = 2057.8462  # H97: 1917.6 -> 2057.8462
$ws.Cells.Item(97,9).Value = 2337.3  # I97: 2115.4167 -> 2337.3
$ws.Cells.Item(97,11).Value = 2337.3  # K97: 2115.4167 -> 2337.3
$ws.Cells.Item(97,13).Value = -1841.3  # M97: -1619.4167 -> -1841.3

# Row 122
$ws.Cells.Item(122,8).Value = 1346.75  # H122: 996.6667 -> 1346.75
$ws.Cells.Item(122,9).Value = 1037  # I122: 880.5 -> 1037
$ws.Cells.Item(122,10).Value = 1863  # J122: 1229 -> 1863
$ws.Cells.Item(122,11).Value = 3111  # K122: 2641.5 -> 3111
$ws.Cells.Item(122,12).Value = 5589  # L122: 3687 -> 5589
$ws.Cells.Item(122,13).Value = -661  # M122: -191.5 -> -661
$ws.Cells.Item(122,14).Value = -10489  # N122: -8587 -> -10489

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40,8).Value = 3760.6155  # H40: 3713.926 -> 3760.6155
$ws.Cells.Item(40,10).Value = 4645.3  # J40: 4450.273 -> 4645.3
$ws.Cells.Item(40,12).Value = 4645.3  # L40: 4450.273 -> 4645.3
$ws.Cells.Item(40,14).Value = -4917.3  # N40: -4722.273 -> -4917.3

# Row 61
$ws.Cells.Item(61,8).Value = 0  # H61: 299.5 -> 0
$ws.Cells.Item(61,9).Value = 0  # I61: 299 -> 0
$ws.Cells.Item(61,10).Value = 0  # J61: 300 -> 0
$ws.Cells.Item(61,11).Value = 0  # K61: 299 -> 0
$ws.Cells.Item(61,12).Value = 0  # L61: 300 -> 0
$ws.Cells.Item(61,13).ClearContents()  # M61 removed (was -97)
$ws.Cells.Item(61,14).ClearContents()  # N61 removed (was -704)

# Row 68
$ws.Cells.Item(68,8).Value = 1566.6666  # H68: 1648.75 -> 1566.6666
$ws.Cells.Item(68,9).Value = 1566.6666  # I68: 1648.75 -> 1566.6666
$ws.Cells.Item(68,11).Value = 1566.6666  # K68: 1648.75 -> 1566.6666
$ws.Cells.Item(68,13).Value = -817.6666  # M68: -899.75 -> -817.6666

# Row 71
$ws.Cells.Item(71,8).Value = 1566.6666  # H71: 1648.75 -> 1566.6666
$ws.Cells.Item(71,9).Value = 1566.6666  # I71: 1648.75 -> 1566.6666
$ws.Cells.Item(71,11).Value = 7833.333000000001  # K71: 8243.75 -> 7833.333000000001
$ws.Cells.Item(71,13).Value = -4089.333000000001  # M71: -4499.75 -> -4089.333000000001

# Row 94
$ws.Cells.Item(94,8).Value = 45000  # H94: 0 -> 45000
$ws.Cells.Item(94,10).Value = 45000  # J94: 0 -> 45000
$ws.Cells.Item(94,12).Value = 45000  # L94: 0 -> 45000
$ws.Cells.Item(94,14).Value = -46352  # N94: None -> -46352

# Row 113
$ws.Cells.Item(113,8).Value = 0  # H113: 299.5 -> 0
$ws.Cells.Item(113,9).Value = 0  # I113: 299 -> 0
$ws.Cells.Item(113,10).Value = 0  # J113: 300 -> 0
$ws.Cells.Item(113,11).Value = 0  # K113: 299 -> 0
$ws.Cells.Item(113,12).Value = 0  # L113: 300 -> 0
$ws.Cells.Item(113,13).ClearContents()  # M113 removed (was 1871)
$ws.Cells.Item(113,14).ClearContents()  # N113 removed (was -4640)

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Cells.Item(4,8).Value = 14341945  # H4: 92857.14 -> 14341945
$ws.Cells.Item(4,9).Value = 16754250  # I4: 225000 -> 16754250
$ws.Cells.Item(4,10).Value = 12532716  # J4: 40000 -> 12532716
$ws.Cells.Item(4,11).Value = 16754250  # K4: 225000 -> 16754250
$ws.Cells.Item(4,12).Value = 12532716  # L4: 40000 -> 12532716
$ws.Cells.Item(4,13).Value = -16754137  # M4: -224887 -> -16754137
$ws.Cells.Item(4,14).Value = -12532942  # N4: -40226 -> -12532942

# Row 41
$ws.Cells.Item(41,8).Value = 24977  # H41: 23488.5 -> 24977
$ws.Cells.Item(41,10).Value = 24977  # J41: 23488.5 -> 24977
$ws.Cells.Item(41,12).Value = 24977  # L41: 23488.5 -> 24977
$ws.Cells.Item(41,14).Value = -25757  # N41: -24268.5 -> -25757

# Row 114
$ws.Cells.Item(114,8).Value = 29997.5  # H114: 48750 -> 29997.5
$ws.Cells.Item(114,10).Value = 29997.5  # J114: 48750 -> 29997.5
$ws.Cells.Item(114,12).Value = 29997.5  # L114: 48750 -> 29997.5
$ws.Cells.Item(114,14).Value = -38675.5  # N114: -57428 -> -38675.5

# Row 122
$ws.Cells.Item(122,8).Value = 9407.637000000001  # H122: 9325.817999999999 -> 9407.637000000001
$ws.Cells.Item(122,9).Value = 3319.7144  # I122: 3545.0667 -> 3319.7144
$ws.Cells.Item(122,10).Value = 20061.5  # J122: 21713.143 -> 20061.5
$ws.Cells.Item(122,11).Value = 9959.143199999999  # K122: 10635.2001 -> 9959.143199999999
$ws.Cells.Item(122,12).Value = 60184.5  # L122: 65139.429 -> 60184.5
$ws.Cells.Item(122,13).Value = -7509.143199999999  # M122: -8185.2001 -> -7509.143199999999
$ws.Cells.Item(122,14).Value = -65084.5  # N122: -70039.429 -> -65084.5
